$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.146.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.532.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.81%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.27"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.534.91"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.87%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -2.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.125"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("E11").Value = "  -5.45%  "

$ws.Range("E12").Value = "  +2.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.136.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000185"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.96%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.538.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.67%  "

$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.276.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.96"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.99%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("E23").Value = "  +3.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.675.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  +7.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.77%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.547.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +6.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.93"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.23%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "168.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.63%  "

$ws.Range("E39").Value = "  +4.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.97"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0798"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.53%  "

$ws.Range("E42").Value = "  -0.25%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.40%  "

$ws.Range("E44").Value = "  -2.61%  "

$ws.Range("E45").Value = "  -0.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.01%  "

$ws.Range("E49").Value = "  +3.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.391.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "302.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.08%  "
